$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 150
$ws.Range("I8").Value = 150
$ws.Range("K8").Value = 450
$ws.Range("M8").Value = -311
$ws.Range("H28").Value = 267.35294
$ws.Range("I28").Value = 259.0625
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 259.0625
$ws.Range("L28").Value = 400
$ws.Range("M28").Value = 225.9375
$ws.Range("N28").Value = -1370
$ws.Range("H33").Value = 1111.0344
$ws.Range("I33").Value = 793.11536
$ws.Range("J33").Value = 3866.3333
$ws.Range("K33").Value = 793.11536
$ws.Range("L33").Value = 3866.3333
$ws.Range("M33").Value = -564.11536
$ws.Range("N33").Value = -4324.3333
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 12000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -12238
$ws.Range("H60").Value = 4000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 4000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 12000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -12968
$ws.Range("H100").Value = 6899.9375
$ws.Range("I100").Value = 2033.3334
$ws.Range("J100").Value = 8023
$ws.Range("K100").Value = 2033.3334
$ws.Range("L100").Value = 8023
$ws.Range("M100").Value = -1492.3334
$ws.Range("N100").Value = -9105
$ws.Range("H132").Value = 4442.517
$ws.Range("I132").Value = 4583.3213
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 13749.9639
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -11219.9639
$ws.Range("N132").Value = -6560
$ws.Range("H137").Value = 1562.5454
$ws.Range("I137").Value = 1324.7646
$ws.Range("J137").Value = 2371
$ws.Range("K137").Value = 3974.2938
$ws.Range("L137").Value = 7113
$ws.Range("M137").Value = -1424.2938
$ws.Range("N137").Value = -12213

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 873.5
$ws.Range("I2").Value = 762.5833
$ws.Range("J2").Value = 1206.25
$ws.Range("K2").Value = 762.5833
$ws.Range("L2").Value = 1206.25
$ws.Range("M2").Value = -649.5833
$ws.Range("N2").Value = -1432.25
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 12
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 9682.493
$ws.Range("I32").Value = 2465.796
$ws.Range("K32").Value = 2465.796
$ws.Range("M32").Value = -2178.796
$ws.Range("H74").Value = 6454305
$ws.Range("I74").Value = 7694613
$ws.Range("J74").Value = 4704.4
$ws.Range("K74").Value = 7694613
$ws.Range("L74").Value = 4704.4
$ws.Range("M74").Value = -7693739
$ws.Range("N74").Value = -6452.4
$ws.Range("H77").Value = 6454305
$ws.Range("I77").Value = 7694613
$ws.Range("J77").Value = 4704.4
$ws.Range("K77").Value = 38473065
$ws.Range("L77").Value = 23522
$ws.Range("M77").Value = -38468697
$ws.Range("N77").Value = -32258
$ws.Range("H102").Value = 2239.5
$ws.Range("I102").Value = 1942.5
$ws.Range("J102").Value = 3130.5
$ws.Range("K102").Value = 1942.5
$ws.Range("L102").Value = 3130.5
$ws.Range("M102").Value = -320.5
$ws.Range("N102").Value = -6374.5
$ws.Range("H116").Value = 873.5
$ws.Range("I116").Value = 762.5833
$ws.Range("J116").Value = 1206.25
$ws.Range("K116").Value = 762.5833
$ws.Range("L116").Value = 1206.25
$ws.Range("M116").Value = 1531.4167
$ws.Range("N116").Value = -5794.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 873.5
$ws.Range("I3").Value = 762.5833
$ws.Range("J3").Value = 1206.25
$ws.Range("K3").Value = 762.5833
$ws.Range("L3").Value = 1206.25
$ws.Range("M3").Value = -648.5833
$ws.Range("N3").Value = -1434.25
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 108.53846
$ws.Range("I22").Value = 91.09999999999999
$ws.Range("J22").Value = 166.66667
$ws.Range("K22").Value = 91.09999999999999
$ws.Range("L22").Value = 166.66667
$ws.Range("M22").Value = 81.90000000000001
$ws.Range("N22").Value = -512.6666700000001
$ws.Range("H80").Value = 107.85714
$ws.Range("I80").Value = 100
$ws.Range("J80").Value = 109.166664
$ws.Range("K80").Value = 100
$ws.Range("L80").Value = 109.166664
$ws.Range("M80").Value = 898
$ws.Range("N80").Value = -2105.166664
$ws.Range("H83").Value = 107.85714
$ws.Range("I83").Value = 100
$ws.Range("J83").Value = 109.166664
$ws.Range("K83").Value = 500
$ws.Range("L83").Value = 545.83332
$ws.Range("M83").Value = 4492
$ws.Range("N83").Value = -10529.83332
$ws.Range("H94").Value = 7479.1377
$ws.Range("I94").Value = 343.6087
$ws.Range("J94").Value = 34832
$ws.Range("K94").Value = 343.6087
$ws.Range("L94").Value = 34832
$ws.Range("M94").Value = 107.3913
$ws.Range("N94").Value = -35734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2386867.2
$ws.Range("I99").Value = 5960282.5
$ws.Range("J99").Value = 4590.3335
$ws.Range("K99").Value = 5960282.5
$ws.Range("L99").Value = 4590.3335
$ws.Range("M99").Value = -5958784.5
$ws.Range("N99").Value = -7586.3335
$ws.Range("H126").Value = 2386867.2
$ws.Range("I126").Value = 5960282.5
$ws.Range("J126").Value = 4590.3335
$ws.Range("K126").Value = 17880847.5
$ws.Range("L126").Value = 13771.0005
$ws.Range("M126").Value = -17878377.5
$ws.Range("N126").Value = -18711.0005
$ws.Range("H132").Value = 2331.9707
$ws.Range("I132").Value = 1408.174
$ws.Range("J132").Value = 4263.5454
$ws.Range("K132").Value = 4224.522
$ws.Range("L132").Value = 12790.6362
$ws.Range("M132").Value = -1694.522
$ws.Range("N132").Value = -17850.6362
$ws.Range("H134").Value = 5317.5884
$ws.Range("I134").Value = 4884.643
$ws.Range("J134").Value = 7338
$ws.Range("K134").Value = 14653.929
$ws.Range("L134").Value = 22014
$ws.Range("M134").Value = -12118.929
$ws.Range("N134").Value = -27084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 100000140
$ws.Range("I6").Value = 111111160
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 333333480
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -333333367
$ws.Range("N6").Value = -3226
$ws.Range("H107").Value = 502.375
$ws.Range("J107").Value = 1028
$ws.Range("L107").Value = 3084
$ws.Range("N107").Value = -6924
$ws.Range("H122").Value = 839.875
$ws.Range("I122").Value = 435
$ws.Range("J122").Value = 1244.75
$ws.Range("K122").Value = 3915
$ws.Range("L122").Value = 11202.75
$ws.Range("M122").Value = -1465
$ws.Range("N122").Value = -16102.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 7500
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 7500
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 7500
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -7990
$ws.Range("H43").Value = 963.4
$ws.Range("I43").Value = 963.4
$ws.Range("K43").Value = 963.4
$ws.Range("M43").Value = -812.4
$ws.Range("H57").Value = 10585.571
$ws.Range("I57").Value = 6000
$ws.Range("J57").Value = 11349.833
$ws.Range("K57").Value = 6000
$ws.Range("L57").Value = 11349.833
$ws.Range("M57").Value = -5180
$ws.Range("N57").Value = -12989.833
$ws.Range("H58").Value = 7692.3076
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 7272.727
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 7272.727
$ws.Range("M58").Value = -9723
$ws.Range("N58").Value = -7826.727
$ws.Range("H80").Value = 2920.6897
$ws.Range("I80").Value = 3375
$ws.Range("J80").Value = 2848
$ws.Range("K80").Value = 3375
$ws.Range("L80").Value = 2848
$ws.Range("M80").Value = -2377
$ws.Range("N80").Value = -4844
$ws.Range("H83").Value = 2920.6897
$ws.Range("I83").Value = 3375
$ws.Range("J83").Value = 2848
$ws.Range("K83").Value = 16875
$ws.Range("L83").Value = 14240
$ws.Range("M83").Value = -11883
$ws.Range("N83").Value = -24224
$ws.Range("H102").Value = 1928.4117
$ws.Range("I102").Value = 2093.5833
$ws.Range("K102").Value = 2093.5833
$ws.Range("M102").Value = -471.5832999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 586.8889
$ws.Range("I22").Value = 541.3333
$ws.Range("J22").Value = 632.44446
$ws.Range("K22").Value = 541.3333
$ws.Range("L22").Value = 632.44446
$ws.Range("M22").Value = -246.3333
$ws.Range("N22").Value = -1222.44446
$ws.Range("H27").Value = 586.8889
$ws.Range("I27").Value = 541.3333
$ws.Range("J27").Value = 632.44446
$ws.Range("K27").Value = 541.3333
$ws.Range("L27").Value = 632.44446
$ws.Range("M27").Value = -434.3333
$ws.Range("N27").Value = -846.44446
$ws.Range("H93").Value = 1002.9032
$ws.Range("I93").Value = 943.375
$ws.Range("J93").Value = 1207
$ws.Range("K93").Value = 943.375
$ws.Range("L93").Value = 1207
$ws.Range("M93").Value = 304.625
$ws.Range("N93").Value = -3703
$ws.Range("H136").Value = 9624239
$ws.Range("I136").Value = 9624239
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 28872717
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -28870167
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 25001106
$ws.Range("I113").Value = 41667856
$ws.Range("J113").Value = 983.625
$ws.Range("K113").Value = 125003568
$ws.Range("L113").Value = 2950.875
$ws.Range("M113").Value = -125001398
$ws.Range("N113").Value = -7290.875
$ws.Range("H122").Value = 250002130
$ws.Range("I122").Value = 500001000
$ws.Range("J122").Value = 3252.5
$ws.Range("K122").Value = 1500003000
$ws.Range("L122").Value = 9757.5
$ws.Range("M122").Value = -1500000550
$ws.Range("N122").Value = -14657.5
$ws.Range("H135").Value = 35633.332
$ws.Range("J135").Value = 35633.332
$ws.Range("L135").Value = 35633.332
$ws.Range("N135").Value = -45773.332
